$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data values scraped on Thu Jul  6 20:36:12 UTC 2023.
# Cells in column D whose new value could be misread as a plain number are forced
# to remain text (matching the original inlineStr/text storage) by setting a "@"
# (Text) number format immediately before assigning the value.

# Row 2
$ws.Range("D2").Value = '30.291.45'
$ws.Range("E2").Value = '  -0.59%  '
# Row 3
$ws.Range("D3").Value = '1.882.26'
$ws.Range("E3").Value = '  -1.48%  '
# Row 4
$ws.Range("E4").Value = '  +0.04%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.02'
$ws.Range("E5").Value = '  -0.44%  '
# Row 6
$ws.Range("E6").Value = '  +0.02%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4688'
$ws.Range("E7").Value = '  -1.94%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2838'
$ws.Range("E8").Value = '  +0.17%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06605'
$ws.Range("E9").Value = '  -1.46%  '
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.68'
$ws.Range("E10").Value = '  +10.55%  '
# Row 11
$ws.Range("E11").Value = '  +1.54%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '98.28'
$ws.Range("E12").Value = '  -2.78%  '
# Row 13
$ws.Range("D13").Value = '1.882.74'
$ws.Range("E13").Value = '  -1.51%  '
# Row 14
$ws.Range("E14").Value = '  -1.89%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6768'
$ws.Range("E15").Value = '  +1.42%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '285.06'
$ws.Range("E16").Value = '  +10.31%  '
# Row 17
$ws.Range("D17").Value = '30.308.09'
$ws.Range("E17").Value = '  -0.65%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.08%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("E19").Value = '  +0.05%  '
# Row 20
$ws.Range("B20").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C20").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D20").Value = '2.130.59'
$ws.Range("E20").Value = '  -1.25%  '
# Row 21
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.401'
$ws.Range("E21").Value = '  +0.30%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007288'
$ws.Range("E22").Value = '  -2.35%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.03%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.191'
$ws.Range("E24").Value = '  -1.46%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.407'
$ws.Range("E25").Value = '  +0.62%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.64'
$ws.Range("E26").Value = '  +0.54%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.26'
$ws.Range("E27").Value = '  +0.69%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.998'
$ws.Range("E28").Value = '  -2.79%  '
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.371'
$ws.Range("E29").Value = '  -1.13%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09717'
$ws.Range("E30").Value = '  -3.02%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.410'
$ws.Range("E31").Value = '  -8.10%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.486'
$ws.Range("E32").Value = '  -1.46%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.145'
$ws.Range("E33").Value = '  -2.51%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04684'
$ws.Range("E34").Value = '  -0.61%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7078'
$ws.Range("E35").Value = '  -2.29%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.098'
$ws.Range("E36").Value = '  -0.56%  '
# Row 37
$ws.Range("B37").Value = 'Frax'
$ws.Range("C37").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9997'
$ws.Range("E37").Value = '  +0.01%  '
# Row 38
$ws.Range("B38").Value = 'HuobiToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.716'
$ws.Range("E38").Value = '  +0.38%  '
# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01877'
$ws.Range("E39").Value = '  -2.00%  '
# Row 40
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.646'
$ws.Range("E40").Value = '  +6.35%  '
# Row 41
$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.526'
$ws.Range("E41").Value = '  -3.34%  '
# Row 42
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.21'
$ws.Range("E42").Value = '  -3.58%  '
# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.979'
$ws.Range("E43").Value = '  +0.82%  '
# Row 44
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8658'
$ws.Range("E44").Value = '  +0.53%  '
# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  +0.04%  '
# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '103.18'
$ws.Range("E46").Value = '  -2.17%  '
# Row 47
$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4199'
$ws.Range("E47").Value = '  -0.94%  '
# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '992.69'
$ws.Range("E48").Value = '  +7.52%  '
# Row 49
$ws.Range("B49").Value = 'Aptos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.293'
$ws.Range("E49").Value = '  -1.01%  '
# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.224'
$ws.Range("E50").Value = '  +4.83%  '
# Row 51
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.04'
$ws.Range("E51").Value = '  -1.96%  '
